# Auto-generated edit script: updates LevePriceNQ/HQ and related columns
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N16").Value = -5460
$ws.Range("M16").ClearContents()
$ws.Range("I16").Value = 0
$ws.Range("H16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("M21").Value = -33368
$ws.Range("I21").Value = 33836
$ws.Range("K21").Value = 33836
$ws.Range("H21").Value = 33836
$ws.Range("I23").Value = 33836
$ws.Range("H23").Value = 33836
$ws.Range("K23").Value = 33836
$ws.Range("M23").Value = -33602
$ws.Range("H51").Value = 9999.75
$ws.Range("J51").Value = 9999.75
$ws.Range("L51").Value = 9999.75
$ws.Range("N51").Value = -10967.75
$ws.Range("M64").Value = -2227
$ws.Range("I64").Value = 2475
$ws.Range("H64").Value = 2475
$ws.Range("K64").Value = 2475
$ws.Range("K67").Value = 2475
$ws.Range("M67").Value = -1617
$ws.Range("I67").Value = 2475
$ws.Range("H67").Value = 2475

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5835.9165
$ws.Range("K32").Value = 5835.9165
$ws.Range("M32").Value = -5548.9165
$ws.Range("I32").Value = 5835.9165
$ws.Range("L35").Value = 5000
$ws.Range("N35").Value = -5812
$ws.Range("M35").Value = -5099.2
$ws.Range("I35").Value = 5505.2
$ws.Range("H35").Value = 5421
$ws.Range("K35").Value = 5505.2
$ws.Range("J35").Value = 5000
$ws.Range("I46").Value = 12250
$ws.Range("K46").Value = 12250
$ws.Range("H46").Value = 23506.834
$ws.Range("J46").Value = 34763.668
$ws.Range("N46").Value = -35401.668
$ws.Range("L46").Value = 34763.668
$ws.Range("M46").Value = -11931
$ws.Range("H74").Value = 2818.375
$ws.Range("K74").Value = 1670.6
$ws.Range("J74").Value = 4731.3335
$ws.Range("L74").Value = 4731.3335
$ws.Range("N74").Value = -6479.3335
$ws.Range("I74").Value = 1670.6
$ws.Range("M74").Value = -796.5999999999999
$ws.Range("H77").Value = 2818.375
$ws.Range("K77").Value = 8353
$ws.Range("J77").Value = 4731.3335
$ws.Range("L77").Value = 23656.6675
$ws.Range("N77").Value = -32392.6675
$ws.Range("I77").Value = 1670.6
$ws.Range("M77").Value = -3985

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H56").Value = 0
$ws.Range("M86").Value = -2884.923
$ws.Range("I86").Value = 4007.923
$ws.Range("H86").Value = 4099.8887
$ws.Range("K86").Value = 4007.923
$ws.Range("I89").Value = 4007.923
$ws.Range("K89").Value = 20039.615
$ws.Range("H89").Value = 4099.8887
$ws.Range("M89").Value = -14423.615
$ws.Range("H94").Value = 1429
$ws.Range("L94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("M99").Value = -158.25
$ws.Range("I99").Value = 1656.25
$ws.Range("K99").Value = 1656.25
$ws.Range("H99").Value = 1647.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -1223.3334
$ws.Range("H31").Value = 3042.7144
$ws.Range("K31").Value = 1518.3334
$ws.Range("J31").Value = 4186
$ws.Range("L31").Value = 4186
$ws.Range("N31").Value = -4776
$ws.Range("I31").Value = 1518.3334
$ws.Range("L34").Value = 4186
$ws.Range("N34").Value = -4590
$ws.Range("I34").Value = 1518.3334
$ws.Range("M34").Value = -1316.3334
$ws.Range("H34").Value = 3042.7144
$ws.Range("K34").Value = 1518.3334
$ws.Range("J34").Value = 4186
$ws.Range("N54").Value = -38814
$ws.Range("H54").Value = 37498
$ws.Range("J54").Value = 37498
$ws.Range("L54").Value = 37498

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I36").Value = 20
$ws.Range("H36").Value = 20
$ws.Range("K36").Value = 60
$ws.Range("M36").Value = 109
$ws.Range("N121").Value = -10034.2855
$ws.Range("H121").Value = 2366.6667
$ws.Range("J121").Value = 2471.4285
$ws.Range("L121").Value = 7414.2855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 114.55556
$ws.Range("J2").Value = 140.28572
$ws.Range("L2").Value = 140.28572
$ws.Range("N2").Value = -366.28572
$ws.Range("L6").Value = 10000
$ws.Range("N6").Value = -10226
$ws.Range("M6").Value = -1137
$ws.Range("I6").Value = 1250
$ws.Range("H6").Value = 7812.5
$ws.Range("K6").Value = 1250
$ws.Range("J6").Value = 10000
$ws.Range("N16").Value = -10500
$ws.Range("M16").Value = -1000
$ws.Range("I16").Value = 1250
$ws.Range("H16").Value = 7812.5
$ws.Range("K16").Value = 1250
$ws.Range("J16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M31").Value = -930.5
$ws.Range("H31").Value = 1222.5
$ws.Range("K31").Value = 1222.5
$ws.Range("I31").Value = 1222.5
$ws.Range("M37").Value = -945.5
$ws.Range("I37").Value = 1222.5
$ws.Range("H37").Value = 1222.5
$ws.Range("K37").Value = 1222.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 609
$ws.Range("K22").Value = 548.6667
$ws.Range("M22").Value = -253.6667
$ws.Range("I22").Value = 548.6667
$ws.Range("H27").Value = 609
$ws.Range("K27").Value = 548.6667
$ws.Range("I27").Value = 548.6667
$ws.Range("M27").Value = -441.6667
$ws.Range("I46").Value = 529.5
$ws.Range("K46").Value = 529.5
$ws.Range("H46").Value = 471.66666
$ws.Range("M46").Value = -341.5
$ws.Range("I122").Value = 4199.1665
$ws.Range("H122").Value = 4199.1665
$ws.Range("K122").Value = 12597.4995
$ws.Range("M122").Value = -10147.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M126").Value = -15306.2861
$ws.Range("I126").Value = 5925.4287
$ws.Range("H126").Value = 5925.4287
$ws.Range("K126").Value = 17776.2861
$ws.Range("M132").Value = 1307.2
$ws.Range("N132").Value = -8210
$ws.Range("I132").Value = 407.6
$ws.Range("H132").Value = 591.1429000000001
$ws.Range("K132").Value = 1222.8
$ws.Range("J132").Value = 1050
$ws.Range("L132").Value = 3150

